$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Mark "ERD - Diagramm" tasks (row 8) as fully complete (100%) for
# Ruhdorfer Alexander (G8) and Sammer Manuel (AA8)
$ws.Range("G8").Value = 1
$ws.Range("AA8").Value = 1

# Add two new tasks for Lamprecht Daniel in rows 9 and 10 of the
# middle table (M:Q)
$ws.Range("M9").Value = "Erstellung des Logischen Modells"
$ws.Range("N9").Value = 42647
$ws.Range("O9").Value = 0.83333333333333337
$ws.Range("P9").Value = 0.91666666666666663
$ws.Range("Q9").Value = 1

$ws.Range("M10").Value = "Erstellung der Create Tables"
$ws.Range("N10").Value = 42647
$ws.Range("O10").Value = 0.4513888888888889
$ws.Range("P10").Value = 0.48680555555555555
$ws.Range("Q10").Value = 1

# Update the view: scroll and selection position
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("H8").Select()
